$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column B: updated cases query text
$ws.Range("B2").Value2 = 'MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.gender = "FEMALE"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '''') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '''') AS `Trial Code`,
    COALESCE(a.arm_id, '''') AS `Arm`,
    COALESCE(a.arm_drug, '''') AS `Arm Treatment`,
    COALESCE(c.disease, '''') AS `Diagnosis`,
    COALESCE(c.gender, '''') AS `Gender`,
    COALESCE(c.race, '''') AS `Race`,
    COALESCE(c.ethnicity, '''') AS `Ethnicity`'

# Row 2, column C: updated stat query text
$ws.Range("C2").Value2 = 'MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.gender = "FEMALE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials'

# Row height for row 2 changed from 87 to 174
$ws.Rows.Item(2).RowHeight = 174

# Selection changed from C9 to B6
$ws.Range("B6").Select()
